$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44203
$ws.Range("K2").Value = 'Black Amber'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 9000
$ws.Range("O2").Value = 10000
$ws.Range("P2").Value = 9500
$ws.Range("Q2").Value = '$/bandeja 10 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 950
$ws.Range("T2").Value = 10

$ws.Range("D3").Value = 44203
$ws.Range("K3").Value = 'Black Amber'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7500
$ws.Range("P3").Value = 7250
$ws.Range("Q3").Value = '$/caja 10 kilos'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 725
$ws.Range("T3").Value = 10

$ws.Range("D4").Value = 44257
$ws.Range("K4").Value = 'Black Amber'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 68
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9559
$ws.Range("Q4").Value = '$/bandeja 18 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 531
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 44257
$ws.Range("K5").Value = 'Black Amber'
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 444
$ws.Range("T5").Value = 18

$ws.Range("D6").Value = 44235
$ws.Range("K6").Value = 'Lemon'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 13500
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 750
$ws.Range("T6").Value = 18

$ws.Range("D7").Value = 44235
$ws.Range("K7").Value = 'Lemon'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11500
$ws.Range("Q7").Value = '$/bandeja 18 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 639
$ws.Range("T7").Value = 18

$ws.Range("D8").Value = 44223
$ws.Range("K8").Value = 'Black Amber'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 639
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44223
$ws.Range("K9").Value = 'Black Amber'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 556
$ws.Range("T9").Value = 18

$ws.Range("D10").Value = 44266
$ws.Range("K10").Value = 'Black Amber'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 9000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 9500
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 528
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 44258
$ws.Range("K11").Value = 'Black Amber'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 9500
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 528
$ws.Range("T11").Value = 18

$ws.Range("D12").Value = 44202
$ws.Range("K12").Value = 'Black Amber'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 9000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 9500
$ws.Range("Q12").Value = '$/caja 10 kilos'
$ws.Range("R12").Value = 'Provincia de Curicó'
$ws.Range("S12").Value = 950
$ws.Range("T12").Value = 10

$ws.Range("D13").Value = 44265
$ws.Range("K13").Value = 'Black Amber'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 9000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 9500
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 528
$ws.Range("T13").Value = 18

$ws.Range("D14").Value = 44265
$ws.Range("K14").Value = 'Black Amber'
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 8000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 8000
$ws.Range("Q14").Value = '$/bandeja 18 kilos granel'
$ws.Range("R14").Value = 'Provincia de Curicó'
$ws.Range("S14").Value = 444
$ws.Range("T14").Value = 18

$ws.Range("D15").Value = 44246
$ws.Range("K15").Value = 'Angeleno'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = 10500
$ws.Range("O15").Value = 11000
$ws.Range("P15").Value = 10750
$ws.Range("Q15").Value = '$/bandeja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 597
$ws.Range("T15").Value = 18

$ws.Range("D16").Value = 44186
$ws.Range("K16").Value = 'Red Beaut'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 13000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 13417
$ws.Range("Q16").Value = '$/caja 15 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 894
$ws.Range("T16").Value = 15

$ws.Range("D17").Value = 44195
$ws.Range("K17").Value = 'Red Beaut'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 90
$ws.Range("N17").Value = 13500
$ws.Range("O17").Value = 14000
$ws.Range("P17").Value = 13722
$ws.Range("Q17").Value = '$/caja 15 kilos granel'
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 915
$ws.Range("T17").Value = 15

$ws.Range("D18").Value = 44252
$ws.Range("K18").Value = 'Black Amber'
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 55
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 11000
$ws.Range("P18").Value = 10545
$ws.Range("Q18").Value = '$/caja 18 kilos empedrada'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 586
$ws.Range("T18").Value = 18

$ws.Range("D19").Value = 44250
$ws.Range("K19").Value = 'Angeleno'
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 120
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 11000
$ws.Range("P19").Value = 10500
$ws.Range("Q19").Value = '$/bandeja 18 kilos granel'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 583
$ws.Range("T19").Value = 18

$ws.Range("D20").Value = 44188
$ws.Range("K20").Value = 'Red Beaut'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 70
$ws.Range("N20").Value = 12500
$ws.Range("O20").Value = 13000
$ws.Range("P20").Value = 12786
$ws.Range("Q20").Value = '$/caja 15 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 852
$ws.Range("T20").Value = 15

$ws.Range("D21").Value = 44277
$ws.Range("K21").Value = 'Black Amber'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 120
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 9500
$ws.Range("Q21").Value = '$/bandeja 18 kilos granel'
$ws.Range("R21").Value = 'Provincia de Curicó'
$ws.Range("S21").Value = 528
$ws.Range("T21").Value = 18

$ws.Range("D22").Value = 44237
$ws.Range("K22").Value = 'Lemon'
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 30
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 13000
$ws.Range("P22").Value = 12500
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("R22").Value = 'Provincia de Curicó'
$ws.Range("S22").Value = 694
$ws.Range("T22").Value = 18

$ws.Range("D23").Value = 44231
$ws.Range("K23").Value = 'Larry Ann'
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 80
$ws.Range("N23").Value = 6000
$ws.Range("O23").Value = 7000
$ws.Range("P23").Value = 6375
$ws.Range("Q23").Value = '$/bandeja 10 kilos granel'
$ws.Range("R23").Value = 'Región Metropolitana'
$ws.Range("S23").Value = 638
$ws.Range("T23").Value = 10

$ws.Range("D24").Value = 44215
$ws.Range("K24").Value = 'Black Amber'
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 65
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 12462
$ws.Range("Q24").Value = '$/bandeja 18 kilos granel'
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 692
$ws.Range("T24").Value = 18

$ws.Range("D25").Value = 44244
$ws.Range("K25").Value = 'Lemon'
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 60
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 9500
$ws.Range("Q25").Value = '$/caja 16 kilos granel'
$ws.Range("R25").Value = 'Región de O''Higgins'
$ws.Range("S25").Value = 594
$ws.Range("T25").Value = 16

$ws.Range("D26").Value = 44259
$ws.Range("K26").Value = 'Black Amber'
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 80
$ws.Range("N26").Value = 8500
$ws.Range("O26").Value = 9000
$ws.Range("P26").Value = 8750
$ws.Range("Q26").Value = '$/bandeja 18 kilos granel'
$ws.Range("R26").Value = 'Región de O''Higgins'
$ws.Range("S26").Value = 486
$ws.Range("T26").Value = 18

$ws.Range("D27").Value = 44259
$ws.Range("K27").Value = 'Black Amber'
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 40
$ws.Range("N27").Value = 8000
$ws.Range("O27").Value = 8000
$ws.Range("P27").Value = 8000
$ws.Range("Q27").Value = '$/bandeja 18 kilos granel'
$ws.Range("R27").Value = 'Región de O''Higgins'
$ws.Range("S27").Value = 444
$ws.Range("T27").Value = 18

$ws.Range("D28").Value = 44267
$ws.Range("K28").Value = 'Angeleno'
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 9000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 9500
$ws.Range("Q28").Value = '$/bandeja 18 kilos granel'
$ws.Range("R28").Value = 'Región de O''Higgins'
$ws.Range("S28").Value = 528
$ws.Range("T28").Value = 18

$ws.Range("D29").Value = 44267
$ws.Range("K29").Value = 'Angeleno'
$ws.Range("L29").Value = 'Segunda'
$ws.Range("M29").Value = 40
$ws.Range("N29").Value = 8000
$ws.Range("O29").Value = 8000
$ws.Range("P29").Value = 8000
$ws.Range("Q29").Value = '$/bandeja 18 kilos granel'
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("S29").Value = 444
$ws.Range("T29").Value = 18

$ws.Range("D30").Value = 44230
$ws.Range("K30").Value = 'Fortuna'
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 11000
$ws.Range("P30").Value = 10700
$ws.Range("Q30").Value = '$/caja 18 kilos granel'
$ws.Range("R30").Value = 'Región de O''Higgins'
$ws.Range("S30").Value = 594
$ws.Range("T30").Value = 18

$ws.Range("D31").Value = 44211
$ws.Range("K31").Value = 'Black Amber'
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 9500
$ws.Range("O31").Value = 10000
$ws.Range("P31").Value = 9792
$ws.Range("Q31").Value = '$/caja 15 kilos granel'
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 653
$ws.Range("T31").Value = 15

$ws.Range("D32").Value = 44189
$ws.Range("K32").Value = 'Red Beaut'
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 60
$ws.Range("N32").Value = 12000
$ws.Range("O32").Value = 13000
$ws.Range("P32").Value = 12333
$ws.Range("Q32").Value = '$/caja 15 kilos granel'
$ws.Range("R32").Value = 'Región de O''Higgins'
$ws.Range("S32").Value = 822
$ws.Range("T32").Value = 15

$ws.Range("D33").Value = 44224
$ws.Range("K33").Value = 'Black Amber'
$ws.Range("L33").Value = 'Especial'
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 13000
$ws.Range("O33").Value = 14000
$ws.Range("P33").Value = 13500
$ws.Range("Q33").Value = '$/caja 18 kilos granel'
$ws.Range("R33").Value = 'Región de O''Higgins'
$ws.Range("S33").Value = 750
$ws.Range("T33").Value = 18

$ws.Range("D34").Value = 44224
$ws.Range("K34").Value = 'Black Amber'
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 80
$ws.Range("N34").Value = 11000
$ws.Range("O34").Value = 12000
$ws.Range("P34").Value = 11500
$ws.Range("Q34").Value = '$/caja 18 kilos granel'
$ws.Range("R34").Value = 'Región de O''Higgins'
$ws.Range("S34").Value = 639
$ws.Range("T34").Value = 18

$ws.Range("D35").Value = 44224
$ws.Range("K35").Value = 'Black Amber'
$ws.Range("L35").Value = 'Segunda'
$ws.Range("M35").Value = 60
$ws.Range("N35").Value = 10000
$ws.Range("O35").Value = 10000
$ws.Range("P35").Value = 10000
$ws.Range("Q35").Value = '$/caja 18 kilos granel'
$ws.Range("R35").Value = 'Región de O''Higgins'
$ws.Range("S35").Value = 556
$ws.Range("T35").Value = 18

